$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.525.44"
$ws.Range("E2").Value = "  +2.09%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.665.46"
$ws.Range("E3").Value = "  +1.40%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  +0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.26"
$ws.Range("E5").Value = "  +0.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4634"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2568"
$ws.Range("E8").Value = "  -1.28%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06117"
$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.664.20"
$ws.Range("E10").Value = "  +1.30%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06938"
$ws.Range("E11").Value = "  -1.99%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.60"
$ws.Range("E12").Value = "  +0.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.331"
$ws.Range("E13").Value = "  -1.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "74.87"
$ws.Range("E14").Value = "  +1.43%  "

$ws.Range("E15").Value = "  -4.93%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  +0.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  +0.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.545.38"
$ws.Range("E18").Value = "  +2.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006716"
$ws.Range("E19").Value = "  +1.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.35"
$ws.Range("E20").Value = "  +0.56%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.878.07"
$ws.Range("E21").Value = "  +1.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.395"
$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.642"
$ws.Range("E23").Value = "  +0.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.213"
$ws.Range("E24").Value = "  -0.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "134.91"
$ws.Range("E25").Value = "  +1.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.83"
$ws.Range("E26").Value = "  -0.66%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.362"
$ws.Range("E27").Value = "  -1.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.707"
$ws.Range("E28").Value = "  +3.42%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "103.76"
$ws.Range("E29").Value = "  -0.55%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.950"
$ws.Range("E30").Value = "  +2.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07687"
$ws.Range("E31").Value = "  -0.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.584"
$ws.Range("E32").Value = "  +0.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04305"
$ws.Range("E33").Value = "  +0.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9408"
$ws.Range("E35").Value = "  +1.22%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.5976"
$ws.Range("E36").Value = "  +1.36%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9178"
$ws.Range("E37").Value = "  +10.53%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.476"
$ws.Range("E38").Value = "  -3.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "105.42"
$ws.Range("E39").Value = "  +6.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9993"
$ws.Range("E40").Value = "  +0.15%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.833"
$ws.Range("E41").Value = "  +4.48%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01457"
$ws.Range("E42").Value = "  -4.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.007"
$ws.Range("E43").Value = "  +6.78%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3699"
$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1109"
$ws.Range("E45").Value = "  +1.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05253"
$ws.Range("E46").Value = "  +1.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.094"
$ws.Range("E47").Value = "  -0.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "29.87"
$ws.Range("E48").Value = "  +2.29%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.586"
$ws.Range("E49").Value = "  +5.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.003"
$ws.Range("E50").Value = "  +0.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.001"
$ws.Range("E51").Value = "  +0.31%  "
